$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.298.74"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "1.910.60"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "321.65"
$ws.Range("E5").Value = "  -2.87%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.4704"
$ws.Range("E7").Value = "  +2.24%  "
$ws.Range("D8").Value = "0.4044"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").Value = "47.66"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "0.08027"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").Value = "0.9981"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").Value = "22.61"
$ws.Range("E12").Value = "  +4.35%  "
$ws.Range("D13").Value = "1.907.76"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "5.874"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "7.100"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "89.46"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("D18").Value = "0.06631"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "17.61"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "29.308.09"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").Value = "5.519"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "2.202"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "2.183.00"
$ws.Range("E26").Value = "  +4.02%  "
$ws.Range("D27").Value = "154.38"
$ws.Range("D28").Value = "19.77"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "6.027"
$ws.Range("E29").Value = "  +9.87%  "
$ws.Range("D30").Value = "2.100"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "117.82"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").Value = "1.070"
$ws.Range("E32").Value = "  +5.26%  "
$ws.Range("D33").Value = "0.09487"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "1.416"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D35").Value = "3.543"
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("D36").Value = "5.364"
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("D37").Value = "0.06070"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "0.02244"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").Value = "8.178"
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").Value = "0.5833"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").Value = "0.1835"
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").Value = "2.491"
$ws.Range("E43").Value = "  +10.08%  "
$ws.Range("D44").Value = "10.09"
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("D45").Value = "0.07884"
$ws.Range("E45").Value = "  +6.01%  "
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").Value = "0.5494"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "1.916"
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("D50").Value = "113.12"
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("D51").Value = "44.15"
$ws.Range("E51").Value = "  -3.52%  "
